# TC11_Bento_Filter_Diagnosis-PapillaryCarcinoma.xlsx - "Fixed Bento 80 Test scripts"
#
# The three Cypher queries on the "startup" sheet (CasesTab / SamplesTab / FilesTab
# rows, stored in column B) get a trailing "order By ... LIMIT 100" clause appended
# (the FilesTab query's old unsorted "order by f.file_name" line is replaced by the
# new, capitalised/limited version). Row heights grow to fit the extra wrapped line,
# and the active selection moves to the last edited cell (B4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : CasesTab query (column B) ---------------------------------
$casesQuery = $ws.Cells.Item(2, 2).Value()
$ws.Cells.Item(2, 2).Value = $casesQuery + "`n order By ss.study_subject_id ASC LIMIT 100"
$ws.Rows.Item(2).RowHeight = 331.2

# --- Row 3 : SamplesTab query (column B) --------------------------------
$samplesQuery = $ws.Cells.Item(3, 2).Value()
$ws.Cells.Item(3, 2).Value = $samplesQuery + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Rows.Item(3).RowHeight = 360

# --- Row 4 : FilesTab query (column B) ----------------------------------
# Replace the trailing "    order by f.file_name" line with the new clause.
$filesQuery = $ws.Cells.Item(4, 2).Value()
$lastNewline = $filesQuery.LastIndexOf("`n")
$filesQueryStub = $filesQuery.Substring(0, $lastNewline)
$ws.Cells.Item(4, 2).Value = $filesQueryStub + "`n     order By f.file_name ASC LIMIT 100"
$ws.Rows.Item(4).RowHeight = 409.6

# --- Final selection: B4 (last cell touched) ----------------------------
$ws.Range("B4").Select()
